# Updated symbol list on Wed Jan 25 07:38:53 UTC 2023 with GitHub Actions
#
# This script refreshes the Price (column D) and Volume(1h) (column E) figures
# on the crypto exchange-token table (Sheet1) to the latest scraped snapshot.
# The source cells are stored as plain text (e.g. "303.50", "-4.74%") rather
# than numeric/percentage values, so each write forces a Text number format
# before assigning the new string and then restores the cell's original
# ("Normal") style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell address, expected previous text, new text
$updates = @(
    @("D2", "303.50", "303.37"),
    @("E2", "-4.74%", "-4.87%"),
    @("D3", "35.21", "35.41"),
    @("E3", "-2.50%", "-2.00%"),
    @("D4", "5.064", "5.058"),
    @("E4", "-1.98%", "-1.91%"),
    @("D5", "0.08003", "0.08005"),
    @("D6", "1.927", "1.922"),
    @("E6", "-10.29%", "-10.67%"),
    @("D7", "4.049", "4.053"),
    @("E7", "-2.14%", "-1.99%"),
    @("D8", "7.749", "7.747"),
    @("E8", "-3.26%", "-3.29%"),
    @("D9", "2.989", "2.943"),
    @("E9", "6.71%", "5.05%"),
    @("D10", "0.9213", "0.9212"),
    @("E10", "-0.64%", "-0.65%"),
    @("D11", "0.1218", "0.1213"),
    @("E11", "19.48%", "19.32%"),
    @("D12", "0.1847", "0.1844"),
    @("E12", "-2.38%", "-2.54%"),
    @("D13", "0.09624", "0.09429"),
    @("E13", "4.58%", "2.32%"),
    @("D14", "0.03589", "0.03569"),
    @("E14", "-0.92%", "-1.49%"),
    @("D15", "0.09859", "0.09855"),
    @("E15", "-0.66%", "-0.61%"),
    @("D16", "0.001390", "0.001389"),
    @("E16", "-3.88%", "-3.47%"),
    @("D17", "0.005746", "0.005789"),
    @("E17", "0.66%", "0.11%"),
    @("D18", "3.500", "3.489"),
    @("E18", "1.07%", "0.93%"),
    @("E19", "1.10%", "1.05%"),
    @("D20", "0.1283", "0.1282"),
    @("E20", "-1.46%", "-1.47%"),
    @("D21", "5.041", "5.042"),
    @("E21", "-3.19%", "-3.13%"),
    @("D22", "0.2466", "0.2464"),
    @("E22", "12.51%", "12.50%"),
    @("D23", "0.04501", "0.04509"),
    @("E23", "-2.18%", "-1.90%"),
    @("D24", "0.001214", "0.001217"),
    @("E24", "-2.75%", "-2.36%"),
    @("D25", "0.004852", "0.004849"),
    @("E25", "2.54%", "2.44%"),
    @("E26", "-0.10%", "-0.07%"),
    @("E27", "-6.93%", "-6.90%"),
    @("D39", "0.01930", "0.01932"),
    @("E39", "-3.97%", "-3.74%"),
    @("D40", "0.04751", "0.04745"),
    @("E40", "-3.55%", "-3.56%"),
    @("D41", "0.007547", "0.007518"),
    @("E41", "-3.07%", "-3.12%"),
    @("D42", "0.009557", "0.009552"),
    @("E42", "22.28%", "21.93%"),
    @("D43", "0.1332", "0.1329"),
    @("E43", "-4.95%", "-5.07%"),
    @("D44", "0.002109", "0.002110"),
    @("E44", "0.21%", "0.26%"),
    @("D45", "0.01116", "0.01118"),
    @("E45", "-6.55%", "-6.35%"),
    @("D46", "0.00006275", "0.00006286"),
    @("E46", "-2.88%", "-2.84%"),
    @("E47", "-0.09%", "-0.08%"),
    @("E48", "57.92%", "104.35%"),
    @("E49", "-31.43%", "-31.41%"),
    @("D50", "0.00002101", "0.00002100"),
    @("E50", "-0.09%", "-0.08%"),
    @("E51", "-0.09%", "-0.08%")
)

$changed = 0
foreach ($u in $updates) {
    $addr = $u[0]
    $expectedOld = $u[1]
    $newVal = $u[2]

    $rng = $ws.Range($addr)
    $current = $rng.Text()

    if ($current -ne $expectedOld) {
        Write-Output "WARNING: $addr was '$current', expected '$expectedOld' (writing '$newVal' anyway)"
    }

    # Force text storage so numeric-looking / percent-looking strings are not
    # reinterpreted as numbers, then drop back to the default style so no
    # extra formatting is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
    $rng.Style = "Normal"

    $changed = $changed + 1
}

Write-Output "Updated $changed cells on '$($ws.Name)'"
